$d = $word.ActiveDocument

# Replace all occurrences of "customerPeriod" with "yearsRegistered"
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("customerPeriod", $true, $false, $false, $false, $false, $true, 1, $false, "yearsRegistered", 2)
